# Applies the "Updated cryptos list" data refresh:
#  - refreshed Price (D) / Volume(1h) (E) figures for most rows
#  - row 39/40 swapped: Stellar now ranks above Hedera
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'23.791.11"
$ws.Range("E2").Value = "  +2.07%  "

# Row 3
$ws.Range("D3").Value = "'1.655.61"
$ws.Range("E3").Value = "  +2.00%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6
$ws.Range("D6").Value = "'304.60"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("D7").Value = "'0.3822"
$ws.Range("E7").Value = "  +1.82%  "

# Row 8
$ws.Range("D8").Value = "'0.3615"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "'51.26"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  +2.74%  "

# Row 11
$ws.Range("D11").Value = "'0.08222"
$ws.Range("E11").Value = "  +1.07%  "

# Row 12
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").Value = "'22.73"
$ws.Range("E13").Value = "  +2.20%  "

# Row 14
$ws.Range("D14").Value = "'6.542"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15
$ws.Range("D15").Value = "'7.451"
$ws.Range("E15").Value = "  +2.52%  "

# Row 16
$ws.Range("E16").Value = "  +0.50%  "

# Row 17
$ws.Range("D17").Value = "'1.641.00"
$ws.Range("E17").Value = "  +1.43%  "

# Row 18
$ws.Range("D18").Value = "'97.88"
$ws.Range("E18").Value = "  +4.18%  "

# Row 19
$ws.Range("D19").Value = "'0.06973"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'6.779"
$ws.Range("E20").Value = "  +3.89%  "

# Row 21
$ws.Range("D21").Value = "'17.76"
$ws.Range("E21").Value = "  +1.59%  "

# Row 22
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("D23").Value = "'12.73"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").Value = "'23.789.37"
$ws.Range("E24").Value = "  +2.09%  "

# Row 25
$ws.Range("D25").Value = "'2.563"
$ws.Range("E25").Value = "  +3.94%  "

# Row 26
$ws.Range("D26").Value = "'3.086"
$ws.Range("E26").Value = "  +0.46%  "

# Row 27
$ws.Range("D27").Value = "'21.33"
$ws.Range("E27").Value = "  +1.18%  "

# Row 28
$ws.Range("D28").Value = "'150.83"
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("D29").Value = "'5.236"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").Value = "'134.82"
$ws.Range("E30").Value = "  +1.59%  "

# Row 31
$ws.Range("D31").Value = "'1.823.20"
$ws.Range("E31").Value = "  +1.43%  "

# Row 32
$ws.Range("D32").Value = "'6.910"
$ws.Range("E32").Value = "  +2.96%  "

# Row 33
$ws.Range("D33").Value = "'1.088"
$ws.Range("E33").Value = "  +2.76%  "

# Row 34
$ws.Range("D34").Value = "'2.131"
$ws.Range("E34").Value = "  +0.88%  "

# Row 35
$ws.Range("D35").Value = "'11.96"
$ws.Range("E35").Value = "  +6.06%  "

# Row 36
$ws.Range("D36").Value = "'0.02841"
$ws.Range("E36").Value = "  +3.67%  "

# Row 37
$ws.Range("E37").Value = "  +1.93%  "

# Row 38
$ws.Range("D38").Value = "'6.148"
$ws.Range("E38").Value = "  +3.21%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.08841"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.07215"
$ws.Range("E40").Value = "  +1.93%  "

# Row 41
$ws.Range("D41").Value = "'12.97"
$ws.Range("E41").Value = "  +8.07%  "

# Row 42
$ws.Range("D42").Value = "'0.7093"
$ws.Range("E42").Value = "  +1.94%  "

# Row 43
$ws.Range("D43").Value = "'1.341"
$ws.Range("E43").Value = "  +1.18%  "

# Row 44
$ws.Range("D44").Value = "'15.95"
$ws.Range("E44").Value = "  +0.13%  "

# Row 45
$ws.Range("D45").Value = "'0.6568"
$ws.Range("E45").Value = "  +1.99%  "

# Row 46
$ws.Range("D46").Value = "'2.337"
$ws.Range("E46").Value = "  +3.41%  "

# Row 47
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("D48").Value = "'3.968"
$ws.Range("E48").Value = "  +0.47%  "

# Row 49
$ws.Range("D49").Value = "'0.07990"
$ws.Range("E49").Value = "  +0.41%  "

# Row 50
$ws.Range("D50").Value = "'128.81"
$ws.Range("E50").Value = "  +2.45%  "

# Row 51
$ws.Range("E51").Value = "  +1.32%  "
